# Auto-generated edit script: apply Phantom_Profits market-price/profit refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 175.88889
$ws.Range("I11").Value = 175.88889
$ws.Range("K11").Value = 175.88889
$ws.Range("M11").Value = -35.88889

$ws.Range("H43").Value = 5399.3335
$ws.Range("J43").Value = 5399.3335
$ws.Range("L43").Value = 5399.3335
$ws.Range("N43").Value = -5537.3335

$ws.Range("H55").Value = 187
$ws.Range("I55").Value = 205.75
$ws.Range("J55").Value = 149.5
$ws.Range("K55").Value = 205.75
$ws.Range("L55").Value = 149.5
$ws.Range("M55").Value = 8.25
$ws.Range("N55").Value = -577.5

$ws.Range("H80").Value = 3070.5
$ws.Range("I80").Value = 2572.889
$ws.Range("J80").Value = 3710.2856
$ws.Range("K80").Value = 7718.667
$ws.Range("L80").Value = 11130.8568
$ws.Range("M80").Value = -6720.667
$ws.Range("N80").Value = -13126.8568

$ws.Range("H83").Value = 3070.5
$ws.Range("I83").Value = 2572.889
$ws.Range("J83").Value = 3710.2856
$ws.Range("K83").Value = 23156.001
$ws.Range("L83").Value = 33392.5704
$ws.Range("M83").Value = -18164.001
$ws.Range("N83").Value = -43376.5704

$ws.Range("H86").Value = 3882.5715
$ws.Range("I86").Value = 3882.5715
$ws.Range("K86").Value = 3882.5715
$ws.Range("M86").Value = -2759.5715

$ws.Range("H89").Value = 3882.5715
$ws.Range("I89").Value = 3882.5715
$ws.Range("K89").Value = 19412.8575
$ws.Range("M89").Value = -13796.8575

$ws.Range("H111").Value = 4312.375
$ws.Range("I111").Value = 4639.8
$ws.Range("J111").Value = 3766.6667
$ws.Range("K111").Value = 13919.4
$ws.Range("L111").Value = 11300.0001
$ws.Range("M111").Value = -10852.4
$ws.Range("N111").Value = -17434.0001

$ws.Range("H116").Value = 2965.6667
$ws.Range("J116").Value = 2449
$ws.Range("L116").Value = 2449
$ws.Range("N116").Value = -9333

$ws.Range("H138").Value = 2422.2666
$ws.Range("J138").Value = 1986.6945
$ws.Range("L138").Value = 5960.083500000001
$ws.Range("N138").Value = -16240.0835

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1598.5385
$ws.Range("I2").Value = 1336.3334
$ws.Range("K2").Value = 1336.3334
$ws.Range("M2").Value = -1223.3334

$ws.Range("H32").Value = 3677.121
$ws.Range("I32").Value = 3323.3125
$ws.Range("J32").Value = 14999
$ws.Range("K32").Value = 3323.3125
$ws.Range("L32").Value = 14999
$ws.Range("M32").Value = -3036.3125
$ws.Range("N32").Value = -15573

$ws.Range("H98").Value = 27449.5
$ws.Range("J98").Value = 27449.5
$ws.Range("L98").Value = 27449.5
$ws.Range("N98").Value = -33439.5

$ws.Range("H116").Value = 1598.5385
$ws.Range("I116").Value = 1336.3334
$ws.Range("K116").Value = 1336.3334
$ws.Range("M116").Value = 957.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1598.5385
$ws.Range("I3").Value = 1336.3334
$ws.Range("K3").Value = 1336.3334
$ws.Range("M3").Value = -1222.3334

$ws.Range("H20").Value = 1221.8125
$ws.Range("I20").Value = 1470.8182
$ws.Range("J20").Value = 674
$ws.Range("K20").Value = 1470.8182
$ws.Range("L20").Value = 674
$ws.Range("M20").Value = -1223.8182
$ws.Range("N20").Value = -1168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1737.6
$ws.Range("J22").Value = 1641.2
$ws.Range("L22").Value = 1641.2
$ws.Range("N22").Value = -2341.2

$ws.Range("H31").Value = 7303.5454
$ws.Range("I31").Value = 8883.817999999999
$ws.Range("J31").Value = 5723.273
$ws.Range("K31").Value = 8883.817999999999
$ws.Range("L31").Value = 5723.273
$ws.Range("M31").Value = -8588.817999999999
$ws.Range("N31").Value = -6313.273

$ws.Range("H34").Value = 7303.5454
$ws.Range("I34").Value = 8883.817999999999
$ws.Range("J34").Value = 5723.273
$ws.Range("K34").Value = 8883.817999999999
$ws.Range("L34").Value = 5723.273
$ws.Range("M34").Value = -8681.817999999999
$ws.Range("N34").Value = -6127.273

$ws.Range("H121").Value = 29996.455
$ws.Range("J121").Value = 29996.455
$ws.Range("L121").Value = 29996.455
$ws.Range("N121").Value = -32616.455

$ws.Range("H132").Value = 11774250
$ws.Range("I132").Value = 15395392
$ws.Range("J132").Value = 5537.75
$ws.Range("K132").Value = 46186176
$ws.Range("L132").Value = 16613.25
$ws.Range("M132").Value = -46183646
$ws.Range("N132").Value = -21673.25

$ws.Range("H134").Value = 1289.4375
$ws.Range("I134").Value = 1330
$ws.Range("J134").Value = 1005.5
$ws.Range("K134").Value = 3990
$ws.Range("L134").Value = 3016.5
$ws.Range("M134").Value = -1455
$ws.Range("N134").Value = -8086.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 849.5
$ws.Range("I87").Value = 849.5
$ws.Range("K87").Value = 2548.5
$ws.Range("M87").Value = -1300.5

$ws.Range("H90").Value = 849.5
$ws.Range("I90").Value = 849.5
$ws.Range("K90").Value = 7645.5
$ws.Range("M90").Value = -1405.5

$ws.Range("H107").Value = 398.57144
$ws.Range("J107").Value = 411.72223
$ws.Range("L107").Value = 1235.16669
$ws.Range("N107").Value = -5075.16669

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws.Range("H131").Value = 3661.5
$ws.Range("J131").Value = 3548.6667
$ws.Range("L131").Value = 10646.0001
$ws.Range("N131").Value = -20726.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H100").Value = 22953.154
$ws.Range("J100").Value = 23615.916
$ws.Range("L100").Value = 23615.916
$ws.Range("N100").Value = -25779.916

$ws.Range("H102").Value = 3284.7273
$ws.Range("I102").Value = 2613.2
$ws.Range("K102").Value = 2613.2
$ws.Range("M102").Value = -991.1999999999998

$ws.Range("H122").Value = 2609.2307
$ws.Range("I122").Value = 2130
$ws.Range("J122").Value = 3168.3333
$ws.Range("K122").Value = 6390
$ws.Range("L122").Value = 9504.999899999999
$ws.Range("M122").Value = -3940
$ws.Range("N122").Value = -14404.9999

$ws.Range("H128").Value = 106980.2
$ws.Range("J128").Value = 106980.2
$ws.Range("L128").Value = 106980.2
$ws.Range("N128").Value = -116940.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1045.4286
$ws.Range("I68").Value = 1063.6
$ws.Range("K68").Value = 1063.6
$ws.Range("M68").Value = -314.5999999999999

$ws.Range("H71").Value = 1045.4286
$ws.Range("I71").Value = 1063.6
$ws.Range("K71").Value = 5318
$ws.Range("M71").Value = -1574

$ws.Range("H101").Value = 6669.875
$ws.Range("J101").Value = 6669.875
$ws.Range("L101").Value = 6669.875
$ws.Range("N101").Value = -13159.875

$ws.Range("H128").Value = 70255.625
$ws.Range("J128").Value = 70255.625
$ws.Range("L128").Value = 70255.625
$ws.Range("N128").Value = -80215.625

$ws.Range("H132").Value = 2194.9
$ws.Range("I132").Value = 1621.2858
$ws.Range("K132").Value = 4863.857400000001
$ws.Range("M132").Value = -2333.857400000001

$ws.Range("H136").Value = 500000000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 37950
$ws.Range("I82").Value = 10000
$ws.Range("J82").Value = 65900
$ws.Range("K82").Value = 10000
$ws.Range("L82").Value = 65900
$ws.Range("M82").Value = -9617
$ws.Range("N82").Value = -66666

$ws.Range("H85").Value = 37950
$ws.Range("I85").Value = 10000
$ws.Range("J85").Value = 65900
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 65900
$ws.Range("M85").Value = -8674
$ws.Range("N85").Value = -68552

$ws.Range("H124").Value = 26619.75
$ws.Range("J124").Value = 26619.75
$ws.Range("L124").Value = 26619.75
$ws.Range("N124").Value = -36439.75

$ws.Range("H126").Value = 1233.6666
$ws.Range("I126").Value = 1233.6666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3700.9998
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1230.9998
$ws.Range("N126").ClearContents()

$ws.Range("H130").Value = 24329.75
$ws.Range("J130").Value = 24329.75
$ws.Range("L130").Value = 24329.75
$ws.Range("N130").Value = -34369.75

$ws.Range("H136").Value = 18116.334
$ws.Range("I136").Value = 18116.334
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 54349.00199999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -51799.00199999999
$ws.Range("N136").ClearContents()

$ws.Range("H140").Value = 57552.668
$ws.Range("J140").Value = 57552.668
$ws.Range("L140").Value = 57552.668
$ws.Range("N140").Value = -67912.66800000001
